$wb = $excel.ActiveWorkbook

# --- Sheet "Forecast Comparison": update MyForecast (column D) values ---
$wsForecast = $wb.Worksheets.Item("Forecast Comparison")

$wsForecast.Range("D3").Value = 125
$wsForecast.Range("D4").Value = 108
$wsForecast.Range("D5").Value = 101
$wsForecast.Range("D6").Value = 118
$wsForecast.Range("D10").Value = 105
$wsForecast.Range("D11").Value = 89
$wsForecast.Range("D12").Value = 101

# --- Sheet "Summary": update Total Forecast metrics (column B) ---
# These cells hold numeric-looking values stored as text, so a leading
# apostrophe is used to force text entry (matches how Excel keeps a
# numeric-looking value as text) instead of letting it auto-convert to a number.
$wsSummary = $wb.Worksheets.Item("Summary")

$wsSummary.Range("B9").Value = "'1799"
$wsSummary.Range("B10").Value = "'972"
$wsSummary.Range("B11").Value = "'464"
